$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F11").Value = "-"
$ws.Range("F12").Value = "-"
$ws.Range("F14").Value = "-"
$ws.Range("F15").Value = "-"
